{"js": "// Locate the unique tail of the text we want to KEEP:\n//   ...sun.invoke.util.ValueConversions.primitiveConversion(sun.invoke.util.Wrapper, Object, boolean)\" is null\n// Everything that follows in the same run (the duplicated\n// \"java.lang.NullPointerException: ...\" line plus the full Java stack\n// trace) must be removed, while the trailing manual line break\n// (<w:br/>) that follows the run is preserved.\n\nconst body = context.document.body;\nconst anchor =\n  'ValueConversions.primitiveConversion(sun.invoke.util.Wrapper, Object, boolean)\" is null';\n\nconst results = body.search(anchor, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // search() results are returned in document order, so the first hit\n  // is the occurrence on the \"...failed:\" line (the one we keep).\n  const match = results.items[0];\n  const para = match.paragraphs.getFirst();\n\n  const afterMatch = match.getRange(\"After\");\n  const paraEnd = para.getRange(\"End\");\n  const tailRange = afterMatch.expandTo(paraEnd);\n\n  // Split the tail on \"\\n\" so the very last piece (the lone trailing\n  // manual line break character) can be excluded from the deletion.\n  const parts = tailRange.split([\"\\n\"], false, false);\n  parts.load(\"text\");\n  await context.sync();\n\n  const n = parts.items.length;\n  if (n > 1) {\n    const toDelete = parts.items[0].expandTo(parts.items[n - 2]);\n    toDelete.insertText(\"\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the unique tail of the text we want to KEEP:\n#   ...sun.invoke.util.ValueConversions.primitiveConversion(sun.invoke.util.Wrapper, Object, boolean)\" is null\n# Everything in the same run AFTER this point (the duplicated\n# \"java.lang.NullPointerException: ...\" line plus the full Java stack\n# trace) must be removed, while the trailing line break / paragraph mark\n# that follow the run are preserved.\n\n$anchor = 'ValueConversions.primitiveConversion(sun.invoke.util.Wrapper, Object, boolean)\" is null'\n\n$findRange = $d.Content.Duplicate\n$found = $findRange.Find.Execute($anchor)\n\nif ($found) {\n    # Locate the paragraph that actually contains the match so we know\n    # where it really ends (just before its trailing line-break +\n    # paragraph mark), rather than relying on the (possibly narrow)\n    # Paragraphs collection of the found range itself.\n    $paraEnd = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($findRange.Start -ge $p.Range.Start -and $findRange.Start -lt $p.Range.End) {\n            $paraEnd = $p.Range.End\n            break\n        }\n    }\n\n    if ($paraEnd -ne $null) {\n        $cutStart = $findRange.End\n        $cutEnd = $paraEnd - 2\n\n        if ($cutEnd -gt $cutStart) {\n            $delRange = $d.Range($cutStart, $cutEnd)\n            $delRange.Text = \"\"\n        }\n    }\n}\n"}
